# Updated symbol list on Wed Feb 15 22:39:22 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto rows whose values changed in this data pull. The sheet stores
# these figures as plain text (e.g. "314.29", "6.05%"), so each cell is
# forced to Text format before the write and the format is cleared again
# afterwards to avoid leaving stray formatting behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$updates = @(
    @{ Row = 2;  D = "314.29";         E = "6.05%" }
    @{ Row = 3;  D = "44.73";          E = "6.99%" }
    @{ Row = 4;  D = "5.148";          E = "2.38%" }
    @{ Row = 5;  D = "0.08055";        E = "6.54%" }
    @{ Row = 6;  D = "4.525";          E = "3.04%" }
    @{ Row = 7;  D = "1.681";          E = "5.66%" }
    @{ Row = 8;  D = "1.087";          E = "16.96%" }
    @{ Row = 9;  D = $null;            E = "7.25%" }
    @{ Row = 10; D = "0.1916";         E = "4.00%" }
    @{ Row = 11; D = "0.09388";        E = "4.00%" }
    @{ Row = 12; D = "0.04271";        E = "6.39%" }
    @{ Row = 13; D = $null;            E = "-0.94%" }
    @{ Row = 14; D = "0.001312";       E = "2.41%" }
    @{ Row = 15; D = "0.005917";       E = "0.72%" }
    @{ Row = 17; D = $null;            E = "0.86%" }
    @{ Row = 18; D = "2.410";          E = "0.08%" }
    @{ Row = 19; D = "0.3374";         E = "1.64%" }
    @{ Row = 20; D = "8.267";          E = "4.53%" }
    @{ Row = 21; D = "0.1372";         E = "-3.38%" }
    @{ Row = 22; D = "0.3139";         E = "4.57%" }
    @{ Row = 23; D = "0.04196";        E = "3.25%" }
    @{ Row = 24; D = "0.001275";       E = "0.71%" }
    @{ Row = 25; D = "0.004559";       E = "15.09%" }
    @{ Row = 26; D = "0.0001342";      E = "9.03%" }
    @{ Row = 38; D = "0.02698";        E = "11.74%" }
    @{ Row = 39; D = "0.05422";        E = "3.88%" }
    @{ Row = 40; D = "0.005486";       E = "-9.11%" }
    @{ Row = 41; D = "0.007751";       E = "-0.26%" }
    @{ Row = 42; D = "0.1417";         E = "6.38%" }
    @{ Row = 43; D = "0.007335";       E = "-2.71%" }
    @{ Row = 44; D = "0.008590";       E = "9.43%" }
    @{ Row = 45; D = "0.3126";         E = "-2.94%" }
    @{ Row = 46; D = "0.00006786";     E = "0.01%" }
    @{ Row = 47; D = "0.00000000746"; E = "-0.65%" }
    @{ Row = 48; D = "0.06156";        E = "31.74%" }
    @{ Row = 49; D = "0.003977";       E = "-5.39%" }
    @{ Row = 50; D = "0.00002088";     E = "-0.65%" }
    @{ Row = 51; D = "0.0001989";      E = "-0.65%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextValue $u.Row 4 $u.D
    }
    Set-TextValue $u.Row 5 $u.E
}
